$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) values for columns B..V. Column A ("Sexo") stays as-is.
$headers = @{
    "B" = "Estado de la Enfermedad al Momento de la Infeccion por SARS-CoV2"
    "C" = "Neumonia"
    "D" = "Antecedente de Trasplante de CPH"
    "E" = "Neutropenia"
    "F" = "HIPOGAMA"
    "G" = "Dexametasona"
    "H" = "Coinfeccion"
    "I" = "TIPO TRASPLANTE EST AUTOLOGO 1 ALOGENICO 2"
    "J" = "Metilprednisolona"
    "K" = "canula 02 en infusion"
    "L" = "Quimioterapia"
    "M" = "EPOC"
    "N" = "Obesidad"
    "O" = "HTA"
    "P" = "DIABETES"
    "Q" = "UTI"
    "R" = "ARM"
    "S" = "Evolucion"
    "T" = "Viral"
    "U" = "Bacteriana"
    "V" = "late_vs_early"
}

# Row 2 (OR values)
$row2 = @{
    "B" = 1.272238042020146
    "C" = "Inf"
    "D" = 0.2745247695715559
    "E" = 0.1814095767956901
    "F" = 0.5144547528977091
    "G" = 3.997566516229475
    "H" = 0
    "I" = 1
    "J" = 0
    "K" = 8.592200908073556
    "L" = 0.1106897013394808
    "M" = 0
    "N" = 1.211988842045695
    "O" = 0.3711307349059412
    "P" = 2.560109292632742
    "Q" = 4.283007771697678
    "R" = 9.034264144710621
    "S" = 6.32070314364931
    "T" = 0
    "U" = 2.585674341544396
    "V" = "OR"
}

# Row 3 (p values)
$row3 = @{
    "B" = 1
    "C" = 0.09834368530020698
    "D" = 0.2107955310270561
    "E" = 0.09953271926619613
    "F" = 0.6493952271984307
    "G" = 0.1228006076572503
    "H" = 1
    "I" = 1
    "J" = 1
    "K" = 0.03274330327096515
    "L" = 0.06080418437397842
    "M" = 1
    "N" = 0.9999999999999999
    "O" = 0.3864416475972541
    "P" = 0.5760869565217393
    "Q" = 0.1819221967963386
    "R" = 0.06080418437397842
    "S" = 0.142080745341615
    "T" = 1
    "U" = 0.4080966482702923
    "V" = "p"
}

# Insert 10 new blank columns before the current column C ("Estado de la
# Enfermedad...") so the existing C..L content (10 columns) shifts right to
# M..V, matching the diff (old C->new L, old D->new M, ..., old L->new V).
$ws.Columns("C:L").Insert()

foreach ($col in $headers.Keys) {
    $ws.Range($col + "1").Value = $headers[$col]
}
foreach ($col in $row2.Keys) {
    $ws.Range($col + "2").Value = $row2[$col]
}
foreach ($col in $row3.Keys) {
    $ws.Range($col + "3").Value = $row3[$col]
}

Write-Output "done"
